$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 232.5
$ws.Range("I9").Value = 180.5
$ws.Range("J9").Value = 284.5
$ws.Range("K9").Value = 180.5
$ws.Range("L9").Value = 284.5
$ws.Range("M9").Value = -11.5
$ws.Range("N9").Value = -622.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1070
$ws.Range("J17").Value = 1101.919
$ws.Range("L17").Value = 3305.757000000001
$ws.Range("N17").Value = -3641.757000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5684.1904
$ws.Range("I74").Value = 3828.3333
$ws.Range("K74").Value = 3828.3333
$ws.Range("M74").Value = -2892.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5684.1904
$ws.Range("I77").Value = 3828.3333
$ws.Range("K77").Value = 19141.6665
$ws.Range("M77").Value = -14461.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 335.43478
$ws.Range("I107").Value = 345.57895
$ws.Range("J107").Value = 287.25
$ws.Range("K107").Value = 345.57895
$ws.Range("L107").Value = 287.25
$ws.Range("M107").Value = 1574.42105
$ws.Range("N107").Value = -4127.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1113056.4
$ws.Range("I61").Value = 1151230.8
$ws.Range("K61").Value = 1151230.8
$ws.Range("M61").Value = -1151018.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 609737.6
$ws.Range("I132").Value = 742161.75
$ws.Range("J132").Value = 4370.143
$ws.Range("K132").Value = 2226485.25
$ws.Range("L132").Value = 13110.429
$ws.Range("M132").Value = -2223955.25
$ws.Range("N132").Value = -18170.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1113056.4
$ws.Range("I136").Value = 1151230.8
$ws.Range("K136").Value = 3453692.4
$ws.Range("M136").Value = -3451142.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9602.223
$ws.Range("I99").Value = 4576.5835
$ws.Range("J99").Value = 19653.5
$ws.Range("K99").Value = 4576.5835
$ws.Range("L99").Value = 19653.5
$ws.Range("M99").Value = -3078.5835
$ws.Range("N99").Value = -22649.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8335521
$ws.Range("I107").Value = 2077.8438
$ws.Range("J107").Value = 41669292
$ws.Range("K107").Value = 2077.8438
$ws.Range("L107").Value = 41669292
$ws.Range("M107").Value = -157.8438000000001
$ws.Range("N107").Value = -41673132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2108
$ws.Range("I16").Value = 1212
$ws.Range("K16").Value = 1212
$ws.Range("M16").Value = -925

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 921.26086
$ws.Range("I22").Value = 729
$ws.Range("J22").Value = 1466
$ws.Range("K22").Value = 729
$ws.Range("L22").Value = 1466
$ws.Range("M22").Value = -379
$ws.Range("N22").Value = -2166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 108139.42
$ws.Range("I31").Value = 174974.06
$ws.Range("J31").Value = 27937.867
$ws.Range("K31").Value = 174974.06
$ws.Range("L31").Value = 27937.867
$ws.Range("M31").Value = -174679.06
$ws.Range("N31").Value = -28527.867

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 108139.42
$ws.Range("I34").Value = 174974.06
$ws.Range("J34").Value = 27937.867
$ws.Range("K34").Value = 174974.06
$ws.Range("L34").Value = 27937.867
$ws.Range("M34").Value = -174772.06
$ws.Range("N34").Value = -28341.867

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 885308.4
$ws.Range("J58").Value = 6233
$ws.Range("L58").Value = 6233
$ws.Range("N58").Value = -6639

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2108
$ws.Range("I113").Value = 1212
$ws.Range("K113").Value = 1212
$ws.Range("M113").Value = 958

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 21853664
$ws.Range("I132").Value = 25643806
$ws.Range("K132").Value = 76931418
$ws.Range("M132").Value = -76928888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 885308.4
$ws.Range("J136").Value = 6233
$ws.Range("L136").Value = 18699
$ws.Range("N136").Value = -23799

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 481.16666
$ws.Range("I33").Value = 272
$ws.Range("J33").Value = 899.5
$ws.Range("K33").Value = 1632
$ws.Range("L33").Value = 5397
$ws.Range("M33").Value = -1349
$ws.Range("N33").Value = -5963

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8500.933999999999
$ws.Range("I68").Value = 1747.5
$ws.Range("J68").Value = 9539.923000000001
$ws.Range("K68").Value = 5242.5
$ws.Range("L68").Value = 28619.769
$ws.Range("M68").Value = -4431.5
$ws.Range("N68").Value = -30241.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 8500.933999999999
$ws.Range("I71").Value = 1747.5
$ws.Range("J71").Value = 9539.923000000001
$ws.Range("K71").Value = 15727.5
$ws.Range("L71").Value = 85859.307
$ws.Range("M71").Value = -11671.5
$ws.Range("N71").Value = -93971.307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 33999.668
$ws.Range("J52").Value = 43499.5
$ws.Range("L52").Value = 43499.5
$ws.Range("N52").Value = -44017.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 32000
$ws.Range("J58").Value = 32000
$ws.Range("L58").Value = 32000
$ws.Range("N58").Value = -32554

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 271852.16
$ws.Range("I80").Value = 343573.06
$ws.Range("K80").Value = 343573.06
$ws.Range("M80").Value = -342575.06

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 271852.16
$ws.Range("I83").Value = 343573.06
$ws.Range("K83").Value = 1717865.3
$ws.Range("M83").Value = -1712873.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 9980
$ws.Range("I97").Value = 10006.667
$ws.Range("J97").Value = 9900
$ws.Range("K97").Value = 10006.667
$ws.Range("L97").Value = 9900
$ws.Range("M97").Value = -9510.666999999999
$ws.Range("N97").Value = -10892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1753.7646
$ws.Range("I113").Value = 1667.6666
$ws.Range("J113").Value = 2399.5
$ws.Range("K113").Value = 1667.6666
$ws.Range("L113").Value = 2399.5
$ws.Range("M113").Value = 502.3334
$ws.Range("N113").Value = -6739.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 697277.2
$ws.Range("I126").Value = 1192809.9
$ws.Range("J126").Value = 3531.4
$ws.Range("K126").Value = 3578429.7
$ws.Range("L126").Value = 10594.2
$ws.Range("M126").Value = -3575959.7
$ws.Range("N126").Value = -15534.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1666.7778
$ws.Range("I46").Value = 837.5
$ws.Range("J46").Value = 2081.4167
$ws.Range("K46").Value = 837.5
$ws.Range("L46").Value = 2081.4167
$ws.Range("M46").Value = -649.5
$ws.Range("N46").Value = -2457.4167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 18239
$ws.Range("I100").Value = 1886
$ws.Range("J100").Value = 100004
$ws.Range("K100").Value = 1886
$ws.Range("L100").Value = 100004
$ws.Range("M100").Value = -1345

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 647977.4
$ws.Range("I132").Value = 995805.4399999999
$ws.Range("J132").Value = 7241.421
$ws.Range("K132").Value = 2987416.32
$ws.Range("L132").Value = 21724.263
$ws.Range("M132").Value = -2984886.32
$ws.Range("N132").Value = -26784.263

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2857.077
$ws.Range("I107").Value = 1420
$ws.Range("K107").Value = 4260
$ws.Range("M107").Value = -2340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 100426
$ws.Range("J108").Value = 100426
$ws.Range("L108").Value = 100426
$ws.Range("N108").Value = -108106

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6494217.5
$ws.Range("I132").Value = 8750951
$ws.Range("K132").Value = 26252853
$ws.Range("M132").Value = -26250323
